# The deck currently has its slide-master theme ("theme2.xml") using the
# "Integral" color palette, while the notes-master theme ("theme1.xml") uses
# the default "Office Theme" palette. The authored edit swaps the content of
# the two theme parts: the slide-master theme becomes "Office Theme" colors
# and the notes-master theme becomes "Integral" colors (font/format schemes
# in both parts are already identical, so only the 12 theme colors actually
# change).
#
# PowerPoint's object model exposes the slide/master theme's color scheme via
# Slide.ThemeColorScheme (a 12-entry collection: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink, in that order) - editing it through any slide
# updates the single shared theme part used by the presentation's slide
# master. Apply the target "Office Theme" RGB values there.

function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target palette ("Office Theme"), in clrScheme order.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToRgb $officeThemeColors[$i - 1]
}
